$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.120.98"
$ws.Range("E2").Value = "  +2.32%  "

$ws.Range("D3").Value = "3.796.89"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'599.84"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'170.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "3.796.23"
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").Value = "'6.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "'0.453"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "'0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.12%  "

$ws.Range("D14").Value = "'36.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "4.436.57"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").Value = "3.795.41"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "69.140.59"
$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").Value = "'7.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  +4.08%  "

$ws.Range("D22").Value = "'471.61"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'0.708"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.75%  "

$ws.Range("D24").Value = "'84.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").Value = "'10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D30").Value = "3.946.79"
$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").Value = "'7.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("D33").Value = "'2.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").Value = "'30.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("D35").Value = "'9.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.99%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = "3.753.15"
$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("D39").Value = "'3.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.71%  "

$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("D42").Value = "'5.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'1.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("D47").Value = "'44.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.14%  "

$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("D49").Value = "'46.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").Value = "'401.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").Value = "'146.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.45%  "
